$wb = $excel.ActiveWorkbook

# --- "login" sheet: point the login URL at GitHub instead of Facebook, and
# rotate in the new password that goes with the GitHub account. ---
$login = $wb.Worksheets.Item("login")
$login.Range("A2").Value = "https://github.com/login"
$login.Range("C2").Value = 'Xanh$&@Ngoc654'

# --- Add a new "searchrepository" test sheet, cloned from "posttext" so it
# keeps the same layout/drawing, positioned right before "posttext". ---
$post = $wb.Worksheets.Item("posttext")
$postIndex = $post.Index
$post.Copy($post, $null)
$searchRepo = $wb.Worksheets.Item($postIndex)
$searchRepo.Name = "searchrepository"

$searchRepo.Range("B1").Value = "repository name"
$searchRepo.Range("A2").Value = 1
$searchRepo.Range("B2").Value = "CrossBrowser"
$searchRepo.Rows(3).Delete()
